$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Q&A pair 1: "How to load a table using external data?"
$answer1 = "To load a table using external data, follow these steps:`n1. On the Definetab, select AutoLoad and then copy the values from your external data source (e.g., Microsoft Excel) into the clipboard.`n2. In the GEOHometab, select the table from the Select Definition box.`n3. If the table has not been created, select Create Table from Template, pick a Category, and select a Table Template.`n4. Click Clipboard in the Load Rows From section to import all the information from the clipboard into the selected table.`nAlternatively, you can also load data from a .txt file:`n1. Ensure the columns are in sequence with that of the table in GEO and that the depth values are appropriate.`n2. Copy the data to clipboard.`n3. In the GEOHometab, select the table from the Select Definition box.`n4. Click Clipboard from the Load Rows from section.`n5. A confirmation dialog box will open; select `"Tables`" to import the qualitative data.`nNote: The .txt file must be tab delimited and have a depth range that is within the global depth range of the receiving odf."

# New Q&A pair 2: "What are the different types of limit data in the GEO limits?"
$answer2 = "The different types of limit data in the GEO limits are:`n1. Histogram or Differential`n2. Qualitative track (also known as a list)`nThese types of limit data provide visual indications and build information into a database, making it easily accessible for export."

$ws.Range("A105").Value() = "How to load a table using external data?"
$ws.Range("B105").Value() = $answer1
$ws.Rows.Item(105).AutoFit()

$ws.Range("A106").Value() = "How to load a table using external data?"
$ws.Range("B106").Value() = $answer1
$ws.Rows.Item(106).AutoFit()

$ws.Range("A107").Value() = "What are the different types of limit data in the GEO limits?"
$ws.Range("B107").Value() = $answer2
$ws.Rows.Item(107).AutoFit()

$ws.Range("A108").Value() = "What are the different types of limit data in the GEO limits?"
$ws.Range("B108").Value() = $answer2
$ws.Rows.Item(108).AutoFit()
